## update reaction time dists
##
## Updates the Midpoint ("E" column) values for three Upgrading-section
## reaction-time parameters (Hydrogenation reaction time, Etherification &
## hydrolysis reaction time, Ring-opening & hydrolysis reaction time). The
## dependent Lower/Upper bound formulas in columns G/I recompute
## automatically. Also clears the now-stale helper flags in column Q for
## rows 34-54 (row 38's flag is left in place), and restores the sheet's
## scroll position / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the three reaction-time Midpoint values -----------------------
$ws.Range("E35").Value = 9.4     # Hydrogenation reaction time
$ws.Range("E41").Value = 6.1     # Etherification & hydrolysis reaction time
$ws.Range("E46").Value = 19      # Ring-opening & hydrolysis reaction time

# --- Clear the stale "Q" helper-column cells for rows 34-54 ---------------
# (row 38 keeps its helper formula, everything else in the block is cleared)
$ws.Range("Q34:Q37").ClearContents()
$ws.Range("Q39:Q54").ClearContents()

# --- Restore view state (scroll position + active selection) --------------
$win = $excel.ActiveWindow
$win.ScrollRow = 27
$win.ScrollColumn = 1
$ws.Range("E36").Select() | Out-Null
